$d = $word.ActiveDocument

# --- Paragraph "Programa" (Portuguese) ---------------------------------
# Split the single run into three runs separated by manual line breaks.
$p1 = $d.Paragraphs(14).Range
$p1.Find.Execute(
    "temperatura.Propriedades magnéticas", $true, $false, $false, $false, $false,
    $true, 1, $false, "temperatura.^lPropriedades magnéticas", 2)

$p1 = $d.Paragraphs(14).Range
$p1.Find.Execute(
    "magnetostricção.Propriedades térmicas", $true, $false, $false, $false, $false,
    $true, 1, $false, "magnetostricção.^lPropriedades térmicas", 2)

# --- Paragraph "Programa" (italic PT tail + EN translation) ------------
$p2 = $d.Paragraphs(15).Range
$p2.Find.Execute(
    "magnetostricção.Propriedades térmicas", $true, $false, $false, $false, $false,
    $true, 1, $false, "magnetostricção.^lPropriedades térmicas", 2)

$p2 = $d.Paragraphs(15).Range
$p2.Find.Execute(
    "térmica.Electrical properties", $true, $false, $false, $false, $false,
    $true, 1, $false, "térmica.^l^lElectrical properties", 2)

$p2 = $d.Paragraphs(15).Range
$p2.Find.Execute(
    "temperature.Magnetic properties", $true, $false, $false, $false, $false,
    $true, 1, $false, "temperature.^lMagnetic properties", 2)

$p2 = $d.Paragraphs(15).Range
$p2.Find.Execute(
    "measurements.Thermal properties", $true, $false, $false, $false, $false,
    $true, 1, $false, "measurements.^lThermal properties", 2)

# --- Paragraph "Bibliografia" -------------------------------------------
$p3 = $d.Paragraphs(19).Range
$p3.Find.Execute(
    "2000.RAYMOND", $true, $false, $false, $false, $false,
    $true, 1, $false, "2000.^lRAYMOND", 2)

$p3 = $d.Paragraphs(19).Range
$p3.Find.Execute(
    "2005.SOLYMAR", $true, $false, $false, $false, $false,
    $true, 1, $false, "2005.^lSOLYMAR", 2)

$p3 = $d.Paragraphs(19).Range
$p3.Find.Execute(
    "2009.NICOLA", $true, $false, $false, $false, $false,
    $true, 1, $false, "2009.^lNICOLA", 2)

$p3 = $d.Paragraphs(19).Range
$p3.Find.Execute(
    "2011ROBERT", $true, $false, $false, $false, $false,
    $true, 1, $false, "2011^lROBERT", 2)

$p3 = $d.Paragraphs(19).Range
$p3.Find.Execute(
    "1998.SPEYER", $true, $false, $false, $false, $false,
    $true, 1, $false, "1998.^lSPEYER", 2)
